# "Got actual dimensions working!" — update the barchart/calendar treatment
# rows on the "Web Parameters" sheet so the pixel / inch dimensions reflect
# the actual (smaller, centered) visualization size, and move the on-screen
# selection/zoom to where the author was looking when they made the change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web Parameters")
$ws.Activate()

# --- Row 5 (treatment 2 / barchart) ---
$ws.Range("N5").Value = 100
$ws.Range("O5").Value = 100
$ws.Range("R5").Value = 8
$ws.Range("S5").Value = 8
$ws.Range("T5").Value = 8.5
$ws.Range("U5").Value = 8.5

# --- Row 6 (treatment 2 / barchart) ---
$ws.Range("N6").Value = 100
$ws.Range("O6").Value = 100
$ws.Range("R6").Value = 8
$ws.Range("S6").Value = 8
$ws.Range("T6").Value = 8.5
$ws.Range("U6").Value = 8.5

# --- Row 7 (treatment 2 / barchart) ---
$ws.Range("N7").Value = 100
$ws.Range("O7").Value = 100
$ws.Range("R7").Value = 4
$ws.Range("S7").Value = 4
$ws.Range("T7").Value = 4.5
$ws.Range("U7").Value = 4.5

# --- Row 8 (treatment 3 / calendar) ---
$ws.Range("N8").Value = 100
$ws.Range("O8").Value = 100
$ws.Range("T8").Value = 8.5
$ws.Range("U8").Value = 8.5

# --- Row 9 (treatment 4 / barchart, drag) ---
$ws.Range("N9").Value = 100
$ws.Range("O9").Value = 100
$ws.Range("R9").Value = 8
$ws.Range("S9").Value = 8
$ws.Range("T9").Value = 8.5
$ws.Range("U9").Value = 8.5

# --- Row 11 (treatment 6 / barchart, titration) ---
$ws.Range("N11").Value = 100
$ws.Range("O11").Value = 100
$ws.Range("R11").Value = 8
$ws.Range("S11").Value = 8
$ws.Range("T11").Value = 8.5
$ws.Range("U11").Value = 8.5

# --- View state: scroll/zoom to where the author ended up, and move the
#     selection from J10 to Q7. ---
$win = $excel.ActiveWindow
$win.Zoom = 120
$excel.Goto($ws.Range("D1"), $true)
$ws.Range("Q7").Select()
